# Add two new IoT error-code rows (412422 / 412423) beneath the existing
# table, and move the active selection to the new last entry (B25),
# matching the author's "finished more IoT methods" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 412422
$ws.Range("B24").Value = "NotFoundException - The specified resource does not exist."

$ws.Range("A25").Value = 412423
$ws.Range("B25").Value = "You must specify a version number greater than 0."

$ws.Range("B25").Select() | Out-Null
